$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("forms")
$ws.Range("B2").Value = '{"ident":"address-form","elems":[{"ident":"TITLE-60664","elementType":"TITLE","dataType":"STRING","texts":{"label":{"textType":"LABEL","texts":{"EN":"Address","DE":"Adresse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[{"extraProp":"SIZE","value":"HUGE"}]},"value":"","required":false,"inline":false,"readOnly":true,"layoutWide":"SIXTEEN","elemEntries":{"entries":[]},"validations":{"hasValidations":false,"rules":[]}},{"ident":"TEXTFIELD-42","elementType":"TEXTFIELD","dataType":"STRING","texts":{"label":{"textType":"LABEL","texts":{"EN":"Street","DE":"Strasse"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[]},"value":"Sonnenweg","required":true,"inline":false,"readOnly":true,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"hasValidations":true,"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"TEXTFIELD-14911","elementType":"TEXTFIELD","dataType":"STRING","texts":{"label":{"textType":"LABEL","texts":{"EN":"Number","DE":"Nummer"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[]},"value":"23a","required":false,"inline":false,"readOnly":true,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"hasValidations":true,"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"TEXTFIELD-33671","elementType":"TEXTFIELD","dataType":"STRING","texts":{"label":{"textType":"LABEL","texts":{"EN":"Postcode","DE":"Plz"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[]},"value":"6414","required":false,"inline":false,"readOnly":true,"layoutWide":"FOUR","elemEntries":{"entries":[]},"validations":{"hasValidations":true,"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":true,"params":{"intParam1":1000,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}},{"ident":"TEXTFIELD-90871","elementType":"TEXTFIELD","dataType":"STRING","texts":{"label":{"textType":"LABEL","texts":{"EN":"Town","DE":"Ort"}},"placeholder":{"textType":"PLACEHOLDER","texts":{"EN":"","DE":""}},"tooltip":{"textType":"TOOLTIP","texts":{"EN":"","DE":""}}},"extras":{"propValues":[]},"value":"Oberarth","required":false,"inline":false,"readOnly":true,"layoutWide":"TWELVE","elemEntries":{"entries":[]},"validations":{"hasValidations":true,"rules":[{"validationType":"EMAIL","enabled":false,"params":{}},{"validationType":"INTEGER","enabled":false,"params":{"intParam1":0,"intParam2":100}},{"validationType":"REG_EXP","enabled":false,"params":{"stringParam":""}}]}}]}'
